$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the time_taken (column F) timestamps on the "data" sheet ---
$newTimes = @(
    "2021-10-05 14:35:13.110584",
    "2021-10-05 14:35:13.110592",
    "2021-10-05 14:35:13.110595",
    "2021-10-05 14:35:13.110598",
    "2021-10-05 14:35:13.110601",
    "2021-10-05 14:35:13.110604",
    "2021-10-05 14:35:13.110606",
    "2021-10-05 14:35:13.110609",
    "2021-10-05 14:35:13.110612",
    "2021-10-05 14:35:13.110614",
    "2021-10-05 14:35:13.110617",
    "2021-10-05 14:35:13.110619",
    "2021-10-05 14:35:13.110622",
    "2021-10-05 14:35:13.110624",
    "2021-10-05 14:35:13.110627",
    "2021-10-05 14:35:13.110629",
    "2021-10-05 14:35:13.110632",
    "2021-10-05 14:35:13.110635",
    "2021-10-05 14:35:13.110637",
    "2021-10-05 14:35:13.110640",
    "2021-10-05 14:35:13.110642",
    "2021-10-05 14:35:13.110645",
    "2021-10-05 14:35:13.110648",
    "2021-10-05 14:35:13.110650",
    "2021-10-05 14:35:13.110653",
    "2021-10-05 14:35:13.110656",
    "2021-10-05 14:35:13.110658",
    "2021-10-05 14:35:13.110661",
    "2021-10-05 14:35:13.110663",
    "2021-10-05 14:35:13.110666",
    "2021-10-05 14:35:13.110668",
    "2021-10-05 14:35:13.110671",
    "2021-10-05 14:35:13.110674",
    "2021-10-05 14:35:13.110676",
    "2021-10-05 14:35:13.110679",
    "2021-10-05 14:35:13.110681",
    "2021-10-05 14:35:13.110684",
    "2021-10-05 14:35:13.110686",
    "2021-10-05 14:35:13.110689",
    "2021-10-05 14:35:13.110692",
    "2021-10-05 14:35:13.110694",
    "2021-10-05 14:35:13.110697",
    "2021-10-05 14:35:13.110700",
    "2021-10-05 14:35:13.110702",
    "2021-10-05 14:35:13.110705",
    "2021-10-05 14:35:13.110707",
    "2021-10-05 14:35:13.110710"
)
for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}


# --- Add the "metadata" sheet (positioned after "data") ---
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$newSheet.Name = "metadata"

# Clone header-row formatting (bold, bordered, centered) from the data sheet's header
$dataSheet.Range("B1").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("B1:G1").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Phagocyte Defects"
$newSheet.Range("C2").Value = 233
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "1.1"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "2021-07-28T07:52:37.846704Z"
$newSheet.Range("F2").Value = "2021-10-05 14:35:13.106788"
$newSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/233/?format=json"
